$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 14 ("Socket") updated from a 14 pin part to a 28 pin part (new Sparkfun SKU)
$ws.Range("B14").Value = "28 pin"
$ws.Range("G14").Value = "PRT-09175"
$ws.Range("H14").Value = "https://www.sparkfun.com/products/9175"
$ws.Hyperlinks.Add($ws.Range("G14"), "https://www.sparkfun.com/products/9175") | Out-Null
$ws.Range("G14").Style = "Hyperlink"

# Price/quantity/total changes for the new part
$ws.Range("I14").Value = 2.95
$ws.Range("J14").Value = 1

# Row 15 quantity bumped (total recalculates automatically via the Table3 formula)
$ws.Range("J15").Value = 45

# Move the active selection to A14 (reflects where the edit was made)
$ws.Range("A14").Select()
